$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 505.39
$ws.Range("I15").Value = 505.39
$ws.Range("K15").Value = 1516.17
$ws.Range("M15").Value = -1347.17
$ws.Range("H19").Value = 1185.2667
$ws.Range("I19").Value = 799.8570999999999
$ws.Range("J19").Value = 1522.5
$ws.Range("K19").Value = 799.8570999999999
$ws.Range("L19").Value = 1522.5
$ws.Range("M19").Value = -624.8570999999999
$ws.Range("N19").Value = -1872.5
$ws.Range("H28").Value = 557.75
$ws.Range("I28").Value = 286.33334
$ws.Range("J28").Value = 2186.25
$ws.Range("K28").Value = 286.33334
$ws.Range("L28").Value = 2186.25
$ws.Range("M28").Value = 198.66666
$ws.Range("N28").Value = -3156.25
$ws.Range("H98").Value = 764.8461
$ws.Range("I98").Value = 673.7
$ws.Range("J98").Value = 1068.6666
$ws.Range("K98").Value = 673.7
$ws.Range("L98").Value = 1068.6666
$ws.Range("M98").Value = 824.3
$ws.Range("N98").Value = -4064.6666
$ws.Range("H103").Value = 1126.5
$ws.Range("J103").Value = 325.42856
$ws.Range("L103").Value = 976.28568
$ws.Range("N103").Value = -2148.28568
$ws.Range("H122").Value = 764.8461
$ws.Range("I122").Value = 673.7
$ws.Range("J122").Value = 1068.6666
$ws.Range("K122").Value = 2021.1
$ws.Range("L122").Value = 3205.9998
$ws.Range("M122").Value = 428.8999999999999
$ws.Range("N122").Value = -8105.9998
$ws.Range("H129").Value = 2314.2922
$ws.Range("I129").Value = 5248
$ws.Range("J129").Value = 914.11365
$ws.Range("K129").Value = 15744
$ws.Range("L129").Value = 2742.34095
$ws.Range("M129").Value = -10744
$ws.Range("N129").Value = -12742.34095
$ws.Range("H132").Value = 4314703.5
$ws.Range("I132").Value = 5323852
$ws.Range("J132").Value = 2887.7273
$ws.Range("K132").Value = 15971556
$ws.Range("L132").Value = 8663.1819
$ws.Range("M132").Value = -15969026
$ws.Range("N132").Value = -13723.1819
$ws.Range("H137").Value = 1464.1428
$ws.Range("I137").Value = 1678.4166
$ws.Range("J137").Value = 1178.4445
$ws.Range("K137").Value = 5035.2498
$ws.Range("L137").Value = 3535.3335
$ws.Range("M137").Value = -2485.2498
$ws.Range("N137").Value = -8635.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22012.479
$ws.Range("I32").Value = 5789.393
$ws.Range("K32").Value = 5789.393
$ws.Range("M32").Value = -5502.393
$ws.Range("H46").Value = 3582.3333
$ws.Range("H63").Value = 2462.2727
$ws.Range("I63").Value = 1983.5714
$ws.Range("J63").Value = 3300
$ws.Range("K63").Value = 1983.5714
$ws.Range("L63").Value = 3300
$ws.Range("M63").Value = -1297.5714
$ws.Range("N63").Value = -4672
$ws.Range("H66").Value = 2462.2727
$ws.Range("I66").Value = 1983.5714
$ws.Range("J66").Value = 3300
$ws.Range("K66").Value = 9917.857
$ws.Range("L66").Value = 16500
$ws.Range("M66").Value = -6485.857
$ws.Range("N66").Value = -23364
$ws.Range("H74").Value = 885.2174
$ws.Range("I74").Value = 784.375
$ws.Range("K74").Value = 784.375
$ws.Range("M74").Value = 89.625
$ws.Range("H77").Value = 885.2174
$ws.Range("I77").Value = 784.375
$ws.Range("K77").Value = 3921.875
$ws.Range("M77").Value = 446.125
$ws.Range("H110").Value = 100210380
$ws.Range("I110").Value = 100210380
$ws.Range("K110").Value = 100210380
$ws.Range("M110").Value = -100208335
$ws.Range("H132").Value = 3995.359
$ws.Range("I132").Value = 3835.8386
$ws.Range("K132").Value = 11507.5158
$ws.Range("M132").Value = -8977.515800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 362.67856
$ws.Range("I64").Value = 543.5
$ws.Range("J64").Value = 290.35
$ws.Range("K64").Value = 543.5
$ws.Range("L64").Value = 290.35
$ws.Range("M64").Value = -318.5
$ws.Range("N64").Value = -740.35
$ws.Range("H67").Value = 362.67856
$ws.Range("I67").Value = 543.5
$ws.Range("J67").Value = 290.35
$ws.Range("K67").Value = 543.5
$ws.Range("L67").Value = 290.35
$ws.Range("M67").Value = 236.5
$ws.Range("N67").Value = -1850.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 7775.231
$ws.Range("I51").Value = 6000
$ws.Range("J51").Value = 7923.1665
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 7923.1665
$ws.Range("M51").Value = -5264
$ws.Range("N51").Value = -9395.166499999999
$ws.Range("H61").Value = 7775.231
$ws.Range("I61").Value = 6000
$ws.Range("J61").Value = 7923.1665
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 7923.1665
$ws.Range("M61").Value = -5652
$ws.Range("N61").Value = -8619.166499999999
$ws.Range("H105").Value = 1149.1904
$ws.Range("I105").Value = 857.9231
$ws.Range("J105").Value = 1622.5
$ws.Range("K105").Value = 857.9231
$ws.Range("L105").Value = 1622.5
$ws.Range("M105").Value = 889.0769
$ws.Range("N105").Value = -5116.5
$ws.Range("H132").Value = 4200.476
$ws.Range("I132").Value = 4660.8667
$ws.Range("J132").Value = 3049.5
$ws.Range("K132").Value = 13982.6001
$ws.Range("L132").Value = 9148.5
$ws.Range("M132").Value = -11452.6001
$ws.Range("N132").Value = -14208.5
$ws.Range("H134").Value = 926.5
$ws.Range("I134").Value = 487.42856
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 1462.28568
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = 1072.71432
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2388.2856
$ws.Range("I68").Value = 1274.5
$ws.Range("J68").Value = 2833.8
$ws.Range("K68").Value = 3823.5
$ws.Range("L68").Value = 8501.400000000001
$ws.Range("M68").Value = -3012.5
$ws.Range("N68").Value = -10123.4
$ws.Range("H71").Value = 2388.2856
$ws.Range("I71").Value = 1274.5
$ws.Range("J71").Value = 2833.8
$ws.Range("K71").Value = 11470.5
$ws.Range("L71").Value = 25504.2
$ws.Range("M71").Value = -7414.5
$ws.Range("N71").Value = -33616.2
$ws.Range("H131").Value = 856.5599999999999
$ws.Range("J131").Value = 861.1717
$ws.Range("L131").Value = 2583.5151
$ws.Range("N131").Value = -12663.5151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111227200
$ws.Range("I80").Value = 166839790
$ws.Range("K80").Value = 166839790
$ws.Range("M80").Value = -166838792
$ws.Range("H83").Value = 111227200
$ws.Range("I83").Value = 166839790
$ws.Range("K83").Value = 834198950
$ws.Range("M83").Value = -834193958
$ws.Range("H122").Value = 1968.1154
$ws.Range("I122").Value = 1307
$ws.Range("J122").Value = 2869.6365
$ws.Range("K122").Value = 3921
$ws.Range("L122").Value = 8608.9095
$ws.Range("M122").Value = -1471
$ws.Range("N122").Value = -13508.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5805.9287
$ws.Range("J7").Value = 6741.6665
$ws.Range("L7").Value = 6741.6665
$ws.Range("N7").Value = -6965.6665
$ws.Range("H57").Value = 9000
$ws.Range("I57").Value = 4000
$ws.Range("J57").Value = 14000
$ws.Range("K57").Value = 4000
$ws.Range("L57").Value = 14000
$ws.Range("M57").Value = -3434
$ws.Range("N57").Value = -15132
$ws.Range("H126").Value = 5805.9287
$ws.Range("J126").Value = 6741.6665
$ws.Range("L126").Value = 20224.9995
$ws.Range("N126").Value = -25164.9995
